# StagingTemplates/Staging.Framework_Project.xlsx
#
# The file was moved into StagingTemplates/ and, while the workbook was
# open/re-saved, the two "...SourceKey" header labels on Sheet1 were
# renamed to "...BusinessKey". Reproduce that content edit here.
#
# (Excel-internal bookkeeping such as the workbook window's last-used
# size/position or a sheet's VBA CodeName is not part of the document
# model exposed through this COM surface - it is simply echoed back
# unchanged by the host - so there is nothing to drive for those via
# script here; only the actual cell content is changed.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "FrameworkBusinessKey"
$ws.Range("C2").Value = "ProjectBusinessKey"
